# Update "想去人数" (interested-count) figures in column F across sheets,
# matching the refreshed data snapshot (output generated at 456a3b4).
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 413
$wsExhibit.Range("F8").Value = 1859
$wsExhibit.Range("F9").Value = 800
$wsExhibit.Range("F10").Value = 19
$wsExhibit.Range("F11").Value = 21
$wsExhibit.Range("F12").Value = 1563
$wsExhibit.Range("F13").Value = 1563
$wsExhibit.Range("F15").Value = 35
$wsExhibit.Range("F16").Value = 1372
$wsExhibit.Range("F18").Value = 389
$wsExhibit.Range("F21").Value = 126
$wsExhibit.Range("F22").Value = 6828
$wsExhibit.Range("F23").Value = 7308
$wsExhibit.Range("F24").Value = 18
$wsExhibit.Range("F25").Value = 164
$wsExhibit.Range("F28").Value = 226
$wsExhibit.Range("F29").Value = 5
$wsExhibit.Range("F33").Value = 1342
$wsExhibit.Range("F34").Value = 206
$wsExhibit.Range("F36").Value = 642
$wsExhibit.Range("F39").Value = 290
$wsExhibit.Range("F40").Value = 162
$wsExhibit.Range("F41").Value = 167
$wsExhibit.Range("F43").Value = 96

$wsPerform = $wb.Worksheets.Item("演出")
$wsPerform.Range("F5").Value = 44
$wsPerform.Range("F17").Value = 261

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F11").Value = 413
$wsAll.Range("F13").Value = 1859
$wsAll.Range("F14").Value = 800
$wsAll.Range("F15").Value = 19
$wsAll.Range("F16").Value = 21
$wsAll.Range("F17").Value = 1563
$wsAll.Range("F18").Value = 1563
$wsAll.Range("F20").Value = 35
$wsAll.Range("F21").Value = 1372
$wsAll.Range("F23").Value = 389
$wsAll.Range("F25").Value = 126
$wsAll.Range("F26").Value = 44
$wsAll.Range("F27").Value = 6828
$wsAll.Range("F28").Value = 7308
$wsAll.Range("F29").Value = 164
$wsAll.Range("F30").Value = 226
$wsAll.Range("F31").Value = 1342
$wsAll.Range("F32").Value = 206
$wsAll.Range("F38").Value = 642
$wsAll.Range("F43").Value = 290
$wsAll.Range("F44").Value = 167
$wsAll.Range("F49").Value = 261

Write-Host "Updated interested-count values across sheets."
